# Added Scania and Guy articulated trucks
# - Adds 5 new vehicle rows (14-18) to the tracking table
# - Adds a "Done" column (J) marking existing + new vehicles as reviewed
# - Adds a "Notes" column (K) with a note on the last new vehicle

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New vehicle rows (14-18) -------------------------------------------
# Columns: A=Vehicle, B=Intro Year, C=Year Order, D=Type, E=ID (formula),
#          F=Top Speed, G=Capacity(goods), H=Cost (formula), I=Running Cost (formula)

$newVehicles = @(
    @{ Row = 14; Name = "Scania 3 Series"; Year = 1987; Order = 1; Speed = 65; Capacity = 44 },
    @{ Row = 15; Name = "Scania 4 Series"; Year = 1995; Order = 1; Speed = 68; Capacity = 48 },
    @{ Row = 16; Name = "Scania R Series"; Year = 2004; Order = 1; Speed = 72; Capacity = 54 },
    @{ Row = 17; Name = "Guy Big J4T";     Year = 1964; Order = 1; Speed = 54; Capacity = 36 },
    @{ Row = 18; Name = "Guy Big J6";      Year = 1964; Order = 2; Speed = 54; Capacity = 22 }
)

foreach ($veh in $newVehicles) {
    $r = $veh.Row
    $ws.Cells.Item($r, 1).Value2 = $veh.Name
    $ws.Cells.Item($r, 2).Value2 = $veh.Year
    $ws.Cells.Item($r, 3).Value2 = $veh.Order
    $ws.Cells.Item($r, 4).Value2 = "Heavy Goods"
    $ws.Cells.Item($r, 5).Formula = "=IF(B$r > 1900, ((B$r-1900)*10)+400+C$r, ((B$r-1730)*2)+C$r)+VLOOKUP(D$r,'ID Scheme'!`$A`$2:`$B`$4,2)"
    $ws.Cells.Item($r, 6).Value2 = $veh.Speed
    $ws.Cells.Item($r, 7).Value2 = $veh.Capacity
    $ws.Cells.Item($r, 8).Formula = "=SQRT(F$r*G$r)/`$B`$1"
    $ws.Cells.Item($r, 8).NumberFormat = "0"
    $ws.Cells.Item($r, 9).Formula = "=H$r*0.9"
    $ws.Cells.Item($r, 9).NumberFormat = "0"
}

# --- "Notes" column (K) ---------------------------------------------------
$ws.Range("K3").Value2 = "Notes"
$ws.Range("K3").Font.Bold = $true

$ws.Range("K18").Value2 = "170hp"

# --- "Done" column (J) -----------------------------------------------------
$ws.Range("J3").Value2 = "Done"
$ws.Range("J3").Font.Bold = $true

for ($r = 4; $r -le 17; $r++) {
    $cell = $ws.Cells.Item($r, 10)
    $cell.Value2 = "x"
    $cell.NumberFormat = "0"
}
# Row 18 (Guy Big J6) is not yet marked done, but keep the same number format
$ws.Cells.Item(18, 10).NumberFormat = "0"
